$wb = $excel.ActiveWorkbook

# --- Overview sheet: file 2d15fcfe-...md is now ready for handoff ---
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"

# --- zh-cn sheet: same status flip, plus refreshed handoff datetime ---
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B5").Value = "Ready for handoff"
$zh.Range("D2").Value = "2016-02-29 04:42:10"
$zh.Range("D3").Value = "2016-02-29 04:42:10"
$zh.Range("D5").Value = "2016-02-29 04:42:10"
$zh.Range("D6").Value = "2016-02-29 04:42:10"
$zh.Range("D7").Value = "2016-02-29 04:42:10"

# --- de-de sheet: same status flip, plus refreshed handoff datetime ---
$de = $wb.Worksheets.Item("de-de")
$de.Range("B5").Value = "Ready for handoff"
$de.Range("D2").Value = "2016-02-29 04:42:22"
$de.Range("D3").Value = "2016-02-29 04:42:22"
$de.Range("D5").Value = "2016-02-29 04:42:22"
$de.Range("D6").Value = "2016-02-29 04:42:22"
$de.Range("D7").Value = "2016-02-29 04:42:22"
